$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65 to shift old row 65 data down to row 66
$ws.Rows.Item(65).Insert()

$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 45494.9
$ws.Cells.Item(2, 3).Value = 6
$ws.Cells.Item(2, 4).Value = 2025
$ws.Cells.Item(2, 5).Value = "06/2025"

$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = 35152.23
$ws.Cells.Item(3, 3).Value = 6
$ws.Cells.Item(3, 4).Value = 2025
$ws.Cells.Item(3, 5).Value = "06/2025"

$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = 11320.2
$ws.Cells.Item(4, 3).Value = 6
$ws.Cells.Item(4, 4).Value = 2025
$ws.Cells.Item(4, 5).Value = "06/2025"

$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = 12775.09
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 2025
$ws.Cells.Item(5, 5).Value = "06/2025"

$ws.Cells.Item(6, 1).Value = 6
$ws.Cells.Item(6, 2).Value = 11707.47
$ws.Cells.Item(6, 3).Value = 6
$ws.Cells.Item(6, 4).Value = 2025
$ws.Cells.Item(6, 5).Value = "06/2025"

$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = 20527.69
$ws.Cells.Item(7, 3).Value = 5
$ws.Cells.Item(7, 4).Value = 2025
$ws.Cells.Item(7, 5).Value = "05/2025"

$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = 29728.09
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 4).Value = 2025
$ws.Cells.Item(8, 5).Value = "05/2025"

$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 19666.55
$ws.Cells.Item(9, 3).Value = 5
$ws.Cells.Item(9, 4).Value = 2025
$ws.Cells.Item(9, 5).Value = "05/2025"

$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = 21883.41
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 2025
$ws.Cells.Item(10, 5).Value = "05/2025"

$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = 45507.35
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 2025
$ws.Cells.Item(11, 5).Value = "05/2025"

$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = 39733.42
$ws.Cells.Item(12, 3).Value = 5
$ws.Cells.Item(12, 4).Value = 2025
$ws.Cells.Item(12, 5).Value = "05/2025"

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 15526.25
$ws.Cells.Item(13, 3).Value = 5
$ws.Cells.Item(13, 4).Value = 2025
$ws.Cells.Item(13, 5).Value = "05/2025"

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 10155.95
$ws.Cells.Item(14, 3).Value = 5
$ws.Cells.Item(14, 4).Value = 2025
$ws.Cells.Item(14, 5).Value = "05/2025"

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 34422.14
$ws.Cells.Item(15, 3).Value = 5
$ws.Cells.Item(15, 4).Value = 2025
$ws.Cells.Item(15, 5).Value = "05/2025"

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 28325.88
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 2025
$ws.Cells.Item(16, 5).Value = "05/2025"

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 12000.74
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 2025
$ws.Cells.Item(17, 5).Value = "05/2025"

$ws.Cells.Item(18, 1).Value = 19
$ws.Cells.Item(18, 2).Value = 17756.75
$ws.Cells.Item(18, 3).Value = 5
$ws.Cells.Item(18, 4).Value = 2025
$ws.Cells.Item(18, 5).Value = "05/2025"

$ws.Cells.Item(19, 1).Value = 20
$ws.Cells.Item(19, 2).Value = 23596.93
$ws.Cells.Item(19, 3).Value = 5
$ws.Cells.Item(19, 4).Value = 2025
$ws.Cells.Item(19, 5).Value = "05/2025"

$ws.Cells.Item(20, 1).Value = 21
$ws.Cells.Item(20, 2).Value = 9475.469999999999
$ws.Cells.Item(20, 3).Value = 5
$ws.Cells.Item(20, 4).Value = 2025
$ws.Cells.Item(20, 5).Value = "05/2025"

$ws.Cells.Item(21, 1).Value = 22
$ws.Cells.Item(21, 2).Value = 27766.95
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = 2025
$ws.Cells.Item(21, 5).Value = "05/2025"

$ws.Cells.Item(22, 1).Value = 23
$ws.Cells.Item(22, 2).Value = 11645.65
$ws.Cells.Item(22, 3).Value = 5
$ws.Cells.Item(22, 4).Value = 2025
$ws.Cells.Item(22, 5).Value = "05/2025"

$ws.Cells.Item(23, 1).Value = 26
$ws.Cells.Item(23, 2).Value = 28728.78
$ws.Cells.Item(23, 3).Value = 5
$ws.Cells.Item(23, 4).Value = 2025
$ws.Cells.Item(23, 5).Value = "05/2025"

$ws.Cells.Item(24, 1).Value = 27
$ws.Cells.Item(24, 2).Value = 15789.98
$ws.Cells.Item(24, 3).Value = 5
$ws.Cells.Item(24, 4).Value = 2025
$ws.Cells.Item(24, 5).Value = "05/2025"

$ws.Cells.Item(25, 1).Value = 28
$ws.Cells.Item(25, 2).Value = 24463.01
$ws.Cells.Item(25, 3).Value = 5
$ws.Cells.Item(25, 4).Value = 2025
$ws.Cells.Item(25, 5).Value = "05/2025"

$ws.Cells.Item(26, 1).Value = 29
$ws.Cells.Item(26, 2).Value = 30467.22
$ws.Cells.Item(26, 3).Value = 5
$ws.Cells.Item(26, 4).Value = 2025
$ws.Cells.Item(26, 5).Value = "05/2025"

$ws.Cells.Item(27, 1).Value = 30
$ws.Cells.Item(27, 2).Value = 23705.99
$ws.Cells.Item(27, 3).Value = 5
$ws.Cells.Item(27, 4).Value = 2025
$ws.Cells.Item(27, 5).Value = "05/2025"

$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = 31020.03
$ws.Cells.Item(28, 3).Value = 4
$ws.Cells.Item(28, 4).Value = 2025
$ws.Cells.Item(28, 5).Value = "04/2025"

$ws.Cells.Item(29, 1).Value = 2
$ws.Cells.Item(29, 2).Value = 48750.66
$ws.Cells.Item(29, 3).Value = 4
$ws.Cells.Item(29, 4).Value = 2025
$ws.Cells.Item(29, 5).Value = "04/2025"

$ws.Cells.Item(30, 1).Value = 3
$ws.Cells.Item(30, 2).Value = 35050.2
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = 2025
$ws.Cells.Item(30, 5).Value = "04/2025"

$ws.Cells.Item(31, 1).Value = 4
$ws.Cells.Item(31, 2).Value = 7210.620000000001
$ws.Cells.Item(31, 3).Value = 4
$ws.Cells.Item(31, 4).Value = 2025
$ws.Cells.Item(31, 5).Value = "04/2025"

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = 10042.97
$ws.Cells.Item(32, 3).Value = 4
$ws.Cells.Item(32, 4).Value = 2025
$ws.Cells.Item(32, 5).Value = "04/2025"

$ws.Cells.Item(33, 1).Value = 8
$ws.Cells.Item(33, 2).Value = 7630.52
$ws.Cells.Item(33, 3).Value = 4
$ws.Cells.Item(33, 4).Value = 2025
$ws.Cells.Item(33, 5).Value = "04/2025"

$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = 19095.69
$ws.Cells.Item(34, 3).Value = 4
$ws.Cells.Item(34, 4).Value = 2025
$ws.Cells.Item(34, 5).Value = "04/2025"

$ws.Cells.Item(35, 1).Value = 10
$ws.Cells.Item(35, 2).Value = 16847.65
$ws.Cells.Item(35, 3).Value = 4
$ws.Cells.Item(35, 4).Value = 2025
$ws.Cells.Item(35, 5).Value = "04/2025"

$ws.Cells.Item(36, 1).Value = 11
$ws.Cells.Item(36, 2).Value = 6701.92
$ws.Cells.Item(36, 3).Value = 4
$ws.Cells.Item(36, 4).Value = 2025
$ws.Cells.Item(36, 5).Value = "04/2025"

$ws.Cells.Item(37, 1).Value = 14
$ws.Cells.Item(37, 2).Value = 38115.7
$ws.Cells.Item(37, 3).Value = 4
$ws.Cells.Item(37, 4).Value = 2025
$ws.Cells.Item(37, 5).Value = "04/2025"

$ws.Cells.Item(38, 1).Value = 15
$ws.Cells.Item(38, 2).Value = 12545.33
$ws.Cells.Item(38, 3).Value = 4
$ws.Cells.Item(38, 4).Value = 2025
$ws.Cells.Item(38, 5).Value = "04/2025"

$ws.Cells.Item(39, 1).Value = 16
$ws.Cells.Item(39, 2).Value = 14981.31
$ws.Cells.Item(39, 3).Value = 4
$ws.Cells.Item(39, 4).Value = 2025
$ws.Cells.Item(39, 5).Value = "04/2025"

$ws.Cells.Item(40, 1).Value = 17
$ws.Cells.Item(40, 2).Value = 14369.85
$ws.Cells.Item(40, 3).Value = 4
$ws.Cells.Item(40, 4).Value = 2025
$ws.Cells.Item(40, 5).Value = "04/2025"

$ws.Cells.Item(41, 1).Value = 22
$ws.Cells.Item(41, 2).Value = 18901.75
$ws.Cells.Item(41, 3).Value = 4
$ws.Cells.Item(41, 4).Value = 2025
$ws.Cells.Item(41, 5).Value = "04/2025"

$ws.Cells.Item(42, 1).Value = 23
$ws.Cells.Item(42, 2).Value = 22062.13
$ws.Cells.Item(42, 3).Value = 4
$ws.Cells.Item(42, 4).Value = 2025
$ws.Cells.Item(42, 5).Value = "04/2025"

$ws.Cells.Item(43, 1).Value = 24
$ws.Cells.Item(43, 2).Value = 17945.9
$ws.Cells.Item(43, 3).Value = 4
$ws.Cells.Item(43, 4).Value = 2025
$ws.Cells.Item(43, 5).Value = "04/2025"

$ws.Cells.Item(44, 1).Value = 25
$ws.Cells.Item(44, 2).Value = 6150.49
$ws.Cells.Item(44, 3).Value = 4
$ws.Cells.Item(44, 4).Value = 2025
$ws.Cells.Item(44, 5).Value = "04/2025"

$ws.Cells.Item(45, 1).Value = 28
$ws.Cells.Item(45, 2).Value = 31243.48
$ws.Cells.Item(45, 3).Value = 4
$ws.Cells.Item(45, 4).Value = 2025
$ws.Cells.Item(45, 5).Value = "04/2025"

$ws.Cells.Item(46, 1).Value = 29
$ws.Cells.Item(46, 2).Value = 33546.1
$ws.Cells.Item(46, 3).Value = 4
$ws.Cells.Item(46, 4).Value = 2025
$ws.Cells.Item(46, 5).Value = "04/2025"

$ws.Cells.Item(47, 1).Value = 30
$ws.Cells.Item(47, 2).Value = 32403.27
$ws.Cells.Item(47, 3).Value = 4
$ws.Cells.Item(47, 4).Value = 2025
$ws.Cells.Item(47, 5).Value = "04/2025"

$ws.Cells.Item(48, 1).Value = 3
$ws.Cells.Item(48, 2).Value = 21522.57
$ws.Cells.Item(48, 3).Value = 3
$ws.Cells.Item(48, 4).Value = 2025
$ws.Cells.Item(48, 5).Value = "03/2025"

$ws.Cells.Item(49, 1).Value = 6
$ws.Cells.Item(49, 2).Value = 37651.31
$ws.Cells.Item(49, 3).Value = 3
$ws.Cells.Item(49, 4).Value = 2025
$ws.Cells.Item(49, 5).Value = "03/2025"

$ws.Cells.Item(50, 1).Value = 7
$ws.Cells.Item(50, 2).Value = 13042.5
$ws.Cells.Item(50, 3).Value = 3
$ws.Cells.Item(50, 4).Value = 2025
$ws.Cells.Item(50, 5).Value = "03/2025"

$ws.Cells.Item(51, 1).Value = 10
$ws.Cells.Item(51, 2).Value = 20815.25
$ws.Cells.Item(51, 3).Value = 3
$ws.Cells.Item(51, 4).Value = 2025
$ws.Cells.Item(51, 5).Value = "03/2025"

$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = 10330.44
$ws.Cells.Item(52, 3).Value = 3
$ws.Cells.Item(52, 4).Value = 2025
$ws.Cells.Item(52, 5).Value = "03/2025"

$ws.Cells.Item(53, 1).Value = 12
$ws.Cells.Item(53, 2).Value = 25429.94
$ws.Cells.Item(53, 3).Value = 3
$ws.Cells.Item(53, 4).Value = 2025
$ws.Cells.Item(53, 5).Value = "03/2025"

$ws.Cells.Item(54, 1).Value = 13
$ws.Cells.Item(54, 2).Value = 14199.5
$ws.Cells.Item(54, 3).Value = 3
$ws.Cells.Item(54, 4).Value = 2025
$ws.Cells.Item(54, 5).Value = "03/2025"

$ws.Cells.Item(55, 1).Value = 14
$ws.Cells.Item(55, 2).Value = 21396.45
$ws.Cells.Item(55, 3).Value = 3
$ws.Cells.Item(55, 4).Value = 2025
$ws.Cells.Item(55, 5).Value = "03/2025"

$ws.Cells.Item(56, 1).Value = 17
$ws.Cells.Item(56, 2).Value = 26250.86
$ws.Cells.Item(56, 3).Value = 3
$ws.Cells.Item(56, 4).Value = 2025
$ws.Cells.Item(56, 5).Value = "03/2025"

$ws.Cells.Item(57, 1).Value = 18
$ws.Cells.Item(57, 2).Value = 3443.92
$ws.Cells.Item(57, 3).Value = 3
$ws.Cells.Item(57, 4).Value = 2025
$ws.Cells.Item(57, 5).Value = "03/2025"

$ws.Cells.Item(58, 1).Value = 19
$ws.Cells.Item(58, 2).Value = 17409.2
$ws.Cells.Item(58, 3).Value = 3
$ws.Cells.Item(58, 4).Value = 2025
$ws.Cells.Item(58, 5).Value = "03/2025"

$ws.Cells.Item(59, 1).Value = 20
$ws.Cells.Item(59, 2).Value = 5023.65
$ws.Cells.Item(59, 3).Value = 3
$ws.Cells.Item(59, 4).Value = 2025
$ws.Cells.Item(59, 5).Value = "03/2025"

$ws.Cells.Item(60, 1).Value = 21
$ws.Cells.Item(60, 2).Value = 22027.47
$ws.Cells.Item(60, 3).Value = 3
$ws.Cells.Item(60, 4).Value = 2025
$ws.Cells.Item(60, 5).Value = "03/2025"

$ws.Cells.Item(61, 1).Value = 24
$ws.Cells.Item(61, 2).Value = 30601.16
$ws.Cells.Item(61, 3).Value = 3
$ws.Cells.Item(61, 4).Value = 2025
$ws.Cells.Item(61, 5).Value = "03/2025"

$ws.Cells.Item(62, 1).Value = 25
$ws.Cells.Item(62, 2).Value = 27204.51
$ws.Cells.Item(62, 3).Value = 3
$ws.Cells.Item(62, 4).Value = 2025
$ws.Cells.Item(62, 5).Value = "03/2025"

$ws.Cells.Item(63, 1).Value = 26
$ws.Cells.Item(63, 2).Value = 19290.09
$ws.Cells.Item(63, 3).Value = 3
$ws.Cells.Item(63, 4).Value = 2025
$ws.Cells.Item(63, 5).Value = "03/2025"

$ws.Cells.Item(64, 1).Value = 27
$ws.Cells.Item(64, 2).Value = 13162.5
$ws.Cells.Item(64, 3).Value = 3
$ws.Cells.Item(64, 4).Value = 2025
$ws.Cells.Item(64, 5).Value = "03/2025"

$ws.Cells.Item(65, 1).Value = 28
$ws.Cells.Item(65, 2).Value = 11137.4
$ws.Cells.Item(65, 3).Value = 3
$ws.Cells.Item(65, 4).Value = 2025
$ws.Cells.Item(65, 5).Value = "03/2025"

$ws.Cells.Item(66, 1).Value = 31
$ws.Cells.Item(66, 2).Value = 41434.74
$ws.Cells.Item(66, 3).Value = 3
$ws.Cells.Item(66, 4).Value = 2025
$ws.Cells.Item(66, 5).Value = "03/2025"
